$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "siteid" column header to "sitenumber"
$ws.Range("C1").Value = "sitenumber"
